$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Rushing" updates
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Updated stats for S.Michel (row 5)
$rushing.Range("C5").Value = 50
$rushing.Range("D5").Value = 25
$rushing.Range("E5").Value = 11
$rushing.Range("F5").Value = 20

# Row 7 - player name corrected from C.Kupp to M.Brown, RZATT corrected to 0
$rushing.Range("B7").Value = "M.Brown"
$rushing.Range("F7").Value = 0

# Copy the formatting used for the existing numbered rows (bold/border/centered
# style) down onto the two brand-new rows before filling in their data.
$rushing.Range("A7").Copy()
$rushing.Range("A8:A9").PasteSpecial(-4122)  # xlPasteFormats

# New row 8 - M.Sargent
$rushing.Range("A8").Value = 6
$rushing.Range("B8").Value = "M.Sargent"
$rushing.Range("C8").Value = 0
$rushing.Range("D8").Value = 1
$rushing.Range("E8").Value = 1
$rushing.Range("F8").Value = 0

# New row 9 - C.Kupp
$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "C.Kupp"
$rushing.Range("C9").Value = 0
$rushing.Range("D9").Value = 0
$rushing.Range("E9").Value = 1
$rushing.Range("F9").Value = 1

# ---------------------------------------------------------------
# Sheet "Receiving" updates
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# S.Michel (row 3)
$receiving.Range("C3").Value = 17
$receiving.Range("D3").Value = 14

# C.Kupp (row 4)
$receiving.Range("C4").Value = 103
$receiving.Range("D4").Value = 81
$receiving.Range("E4").Value = 36
$receiving.Range("F4").Value = 19
$receiving.Range("G4").Value = 25

# V.Jefferson (row 5)
$receiving.Range("C5").Value = 50
$receiving.Range("D5").Value = 30
$receiving.Range("E5").Value = 22
$receiving.Range("G5").Value = 13
$receiving.Range("H5").Value = 4

# B.Skowronek (row 6)
$receiving.Range("C6").Value = 17
$receiving.Range("E6").Value = 3
$receiving.Range("F6").Value = 2

# O.Beckham (row 7)
$receiving.Range("C7").Value = 49
$receiving.Range("D7").Value = 44
$receiving.Range("E7").Value = 24
$receiving.Range("F7").Value = 10
$receiving.Range("G7").Value = 19
$receiving.Range("H7").Value = 10

# Copy formatting for the new numbered rows (8-10) before filling them in.
$receiving.Range("A7").Copy()
$receiving.Range("A8:A10").PasteSpecial(-4122)  # xlPasteFormats

# New row 8 - K.Blanton (inserted ahead of T.Higbee / J.Mundt)
$receiving.Range("A8").Value = 6
$receiving.Range("B8").Value = "K.Blanton"
$receiving.Range("C8").Value = 2
$receiving.Range("D8").Value = 1
$receiving.Range("E8").Value = 0
$receiving.Range("F8").Value = 0
$receiving.Range("G8").Value = 1
$receiving.Range("H8").Value = 0

# Row 9 now holds T.Higbee (shifted down from row 8)
$receiving.Range("A9").Value = 7
$receiving.Range("B9").Value = "T.Higbee"
$receiving.Range("C9").Value = 54
$receiving.Range("D9").Value = 41
$receiving.Range("E9").Value = 8
$receiving.Range("F9").Value = 3
$receiving.Range("G9").Value = 16
$receiving.Range("H9").Value = 11

# Row 10 now holds J.Mundt (shifted down from row 9)
$receiving.Range("A10").Value = 8
$receiving.Range("B10").Value = "J.Mundt"
$receiving.Range("C10").Value = 1
$receiving.Range("D10").Value = 1
$receiving.Range("E10").Value = 0
$receiving.Range("F10").Value = 0
$receiving.Range("G10").Value = 0
$receiving.Range("H10").Value = 0
